$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2020-Q4" (position 2),
#    which pushes the existing "总计" sheet down to position 3.
#    We clone the "总计" sheet (instead of adding a blank one) so the new
#    sheet correctly inherits the same sheet-level properties (outline /
#    page-setup settings) that the target sheet carries.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(1)
$total = $wb.Worksheets.Item(2)

$total.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. Populate the "2022-Q1" sheet.
#    Extend the already-styled header (column D) and data (columns D:D)
#    formatting rightwards/downwards by tiling a single formatted cell
#    across the new destination ranges, so columns E..H and row 3 pick up
#    the same style indices instead of Excel minting brand-new ones.
# ---------------------------------------------------------------------------
$newSheet.Range("D1").Copy($newSheet.Range("E1:H1"))
$newSheet.Range("A2").Copy($newSheet.Range("A3"))
$newSheet.Range("D2:D3").Copy($newSheet.Range("D2:H3"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B..G hold numeric-looking values that must be stored as TEXT
# (matching the source data), so force text number-formatting on that block
# before assigning the values. Column A (index) and H (rank) stay numeric.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002236"
$newSheet.Range("C2").Value = "大成中证360互联网+大数据100指数A"
$newSheet.Range("D2").Value = "5.67"
$newSheet.Range("E2").Value = "93.32"
$newSheet.Range("F2").Value = "1.00"
$newSheet.Range("G2").Value = "0.0567"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "003359"
$newSheet.Range("C3").Value = "大成中证360互联网+大数据100指数C"
$newSheet.Range("D3").Value = "4.08"
$newSheet.Range("E3").Value = "93.32"
$newSheet.Range("F3").Value = "1.00"
$newSheet.Range("G3").Value = "0.0408"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift the existing "2020-Q4" row down
#    to row 3 (cloning column A's style first), and write a brand-new row 2
#    summarising the "2022-Q1" sheet just added.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(3)
$total.Range("A2").Copy($total.Range("A3"))

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2020-Q4"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 1.34

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.1

# ---------------------------------------------------------------------------
# 4. Restore "2020-Q4" as the active tab (unchanged by the diff).
# ---------------------------------------------------------------------------
$q4.Activate()
